$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Air")
$ws.Range("B3").Value = 42885
